$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.184.30"
$ws.Range("E2").Value = "  -1.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.011.44"
$ws.Range("E3").Value = "  -1.80%  "

$ws.Range("E4").Value = "  -0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.60"
$ws.Range("E5").Value = "  +2.22%  "

$ws.Range("E6").Value = "  -3.45%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.61"
$ws.Range("E7").Value = "  +9.60%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "59.05"
$ws.Range("E9").Value = "  -6.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.368"
$ws.Range("E10").Value = "  -1.05%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0744"
$ws.Range("E11").Value = "  -0.98%  "

$ws.Range("E12").Value = "  -1.78%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.899"
$ws.Range("E13").Value = "  -1.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.79"
$ws.Range("E14").Value = "  +2.23%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.304.70"
$ws.Range("E15").Value = "  -1.97%  "

$ws.Range("E16").Value = "  -0.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.57"
$ws.Range("E17").Value = "  +9.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.038.16"
$ws.Range("E18").Value = "  -0.69%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.075.64"
$ws.Range("E19").Value = "  -1.22%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.94"
$ws.Range("E20").Value = "  +0.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0855"
$ws.Range("E21").Value = "  -0.78%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.25"
$ws.Range("E22").Value = "  +0.46%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.48"
$ws.Range("E23").Value = "  -1.46%  "

$ws.Range("E24").Value = "  +18.48%  "

$ws.Range("E25").Value = "  +0.03%  "

$ws.Range("E26").Value = "  -2.75%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.51"
$ws.Range("E27").Value = "  +2.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.88"
$ws.Range("E28").Value = "  -0.50%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.52"
$ws.Range("E29").Value = "  -2.48%  "

$ws.Range("E30").Value = "  -1.01%  "

$ws.Range("E31").Value = "  +2.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.18"
$ws.Range("E32").Value = "  -1.91%  "

$ws.Range("E33").Value = "  +24.54%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.50"
$ws.Range("E34").Value = "  +1.54%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0602"
$ws.Range("E35").Value = "  +0.24%  "

$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.50"
$ws.Range("E36").Value = "  +12.64%  "

$ws.Range("E37").Value = "  -0.27%  "

$ws.Range("E38").Value = "  -1.15%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.83"
$ws.Range("E39").Value = "  +13.68%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.103"
$ws.Range("E40").Value = "  +12.67%  "

$ws.Range("E41").Value = "  -1.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.90"
$ws.Range("E42").Value = "  +0.53%  "

$ws.Range("E43").Value = "  -1.12%  "

$ws.Range("E44").Value = "  +1.52%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.59"
$ws.Range("E45").Value = "  +3.57%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.427.24"
$ws.Range("E46").Value = "  +3.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "93.76"
$ws.Range("E47").Value = "  -0.47%  "

$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.79"
$ws.Range("E48").Value = "  +4.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.51"
$ws.Range("E49").Value = "  +9.94%  "

$ws.Range("E50").Value = "  -1.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.29"
$ws.Range("E51").Value = "  +2.72%  "
